$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its original text formatting (values like
# "30.221.17" / "1.000" / "0.000007273" must stay text, not be auto-
# converted to numbers). Only the D cells that actually change below are
# touched (D4's price is untouched by this update, so row 4 is skipped).
# Each contiguous run gets its own call: a single multi-area Range
# ("D2:D3,D5:D51") only applied NumberFormat to its first area here.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.221.17"
$ws.Range("E2").Value = "  -0.75%  "
# Row 3
$ws.Range("D3").Value = "1.839.99"
$ws.Range("E3").Value = "  -1.49%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").Value = "232.73"
$ws.Range("E5").Value = "  -1.41%  "
# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.03%  "
# Row 7
$ws.Range("D7").Value = "0.4668"
$ws.Range("E7").Value = "  -3.30%  "
# Row 8
$ws.Range("D8").Value = "0.2710"
$ws.Range("E8").Value = "  -3.24%  "
# Row 9
$ws.Range("D9").Value = "0.06272"
$ws.Range("E9").Value = "  -3.62%  "
# Row 10
$ws.Range("D10").Value = "1.842.20"
$ws.Range("E10").Value = "  -1.19%  "
# Row 11
$ws.Range("D11").Value = "0.07406"
$ws.Range("E11").Value = "  -0.52%  "
# Row 12
$ws.Range("D12").Value = "16.10"
$ws.Range("E12").Value = "  -0.83%  "
# Row 13
$ws.Range("D13").Value = "4.928"
$ws.Range("E13").Value = "  -3.03%  "
# Row 14
$ws.Range("D14").Value = "83.66"
$ws.Range("E14").Value = "  -3.94%  "
# Row 15
$ws.Range("D15").Value = "0.6191"
$ws.Range("E15").Value = "  -3.56%  "
# Row 16
$ws.Range("D16").Value = "30.152.51"
$ws.Range("E16").Value = "  -0.91%  "
# Row 17
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  -0.05%  "
# Row 18
$ws.Range("D18").Value = "226.40"
$ws.Range("E18").Value = "  -1.74%  "
# Row 19
$ws.Range("D19").Value = "0.000007273"
$ws.Range("E19").Value = "  -2.85%  "
# Row 20
$ws.Range("D20").Value = "12.35"
$ws.Range("E20").Value = "  -4.96%  "
# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.082.32"
$ws.Range("E21").Value = "  -0.65%  "
# Row 22
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.03%  "
# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "4.883"
$ws.Range("E23").Value = "  -5.06%  "
# Row 24
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "5.845"
$ws.Range("E24").Value = "  -4.01%  "
# Row 25
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "9.197"
$ws.Range("E25").Value = "  -1.43%  "
# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "164.38"
$ws.Range("E26").Value = "  -3.22%  "
# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "17.74"
$ws.Range("E27").Value = "  -3.15%  "
# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "1.862"
$ws.Range("E28").Value = "  -2.37%  "
# Row 29
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "0.1033"
$ws.Range("E29").Value = "  -0.85%  "
# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.373"
$ws.Range("E30").Value = "  -0.81%  "
# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.077"
$ws.Range("E31").Value = "  -4.75%  "
# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.809"
$ws.Range("E32").Value = "  -4.42%  "
# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.04817"
$ws.Range("E33").Value = "  -3.32%  "
# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.142"
$ws.Range("E34").Value = "  -3.16%  "
# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7079"
$ws.Range("E35").Value = "  -4.53%  "
# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.698"
$ws.Range("E36").Value = "  -0.46%  "
# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.01873"
$ws.Range("E37").Value = "  -3.10%  "
# Row 38
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.650"
$ws.Range("E38").Value = "  +0.65%  "
# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "0.8922"
$ws.Range("E39").Value = "  -2.84%  "
# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "1.921"
$ws.Range("E40").Value = "  -6.19%  "
# Row 41
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "104.15"
$ws.Range("E41").Value = "  -1.58%  "
# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.56%  "
# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.529"
$ws.Range("E43").Value = "  -0.87%  "
# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4009"
$ws.Range("E44").Value = "  -4.42%  "
# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "7.022"
$ws.Range("E45").Value = "  -2.62%  "
# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1194"
$ws.Range("E46").Value = "  -2.74%  "
# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "59.68"
$ws.Range("E47").Value = "  -3.44%  "
# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.520"
$ws.Range("E48").Value = "  -4.10%  "
# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "32.83"
$ws.Range("E49").Value = "  -2.22%  "
# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05512"
$ws.Range("E50").Value = "  -2.41%  "
# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.358"
$ws.Range("E51").Value = "  -4.58%  "

Write-Output "applied 161 cell updates"
